$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'40.260.30"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.37%  '
$ws.Range("D3").Value = "'2.343.99"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -5.21%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'308.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.56%  '
$ws.Range("D6").Value = "'85.82"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.52%  '
$ws.Range("D7").Value = "'0.530"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.47%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = "'0.487"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.33%  '
$ws.Range("D10").Value = "'0.0821"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.62%  '
$ws.Range("D11").Value = "'30.46"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -7.27%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").Value = "'2.708.44"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.22%  '
$ws.Range("D14").Value = "'6.48"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.72%  '
$ws.Range("D15").Value = "'14.81"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.04%  '
$ws.Range("D16").Value = "'2.349.66"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.04%  '
$ws.Range("D17").Value = "'0.758"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.10%  '
$ws.Range("D18").Value = "'40.231.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.33%  '
$ws.Range("D19").Value = "'0.0₃0907"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.39%  '
$ws.Range("D20").Value = "'6.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.78%  '
$ws.Range("D21").Value = "'67.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.89%  '
$ws.Range("D22").Value = "'10.80"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.59%  '
$ws.Range("D23").Value = "'235.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("D24").Value = "'2.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -6.69%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").Value = "'1.82"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.09%  '
$ws.Range("D27").Value = "'23.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -5.32%  '
$ws.Range("D28").Value = "'2.15"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("D29").Value = "'9.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.72%  '
$ws.Range("D30").Value = "'35.09"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("D31").Value = "'152.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.98%  '
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").Value = "'5.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.74%  '
$ws.Range("D34").Value = "'0.0729"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.04%  '
$ws.Range("D35").Value = "'2.45"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.63%  '
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("E37").Value = '  -2.19%  '
$ws.Range("D38").Value = "'15.91"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -6.59%  '
$ws.Range("D39").Value = "'2.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.58%  '
$ws.Range("E40").Value = '  -6.10%  '
$ws.Range("E41").Value = '  -4.20%  '
$ws.Range("E42").Value = '  -5.37%  '
$ws.Range("D43").Value = "'1.954.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").Value = "'0.0268"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.13%  '
$ws.Range("D45").Value = "'17.68"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.86%  '
$ws.Range("E46").Value = '  -0.88%  '
$ws.Range("D47").Value = "'2.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -8.82%  '
$ws.Range("D48").Value = "'2.570.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.99%  '
$ws.Range("D49").Value = "'93.22"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.22%  '
$ws.Range("D50").Value = "'71.60"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.54%  '
$ws.Range("D51").Value = "'50.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.98%  '
